$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selección")
Write-Host $ws.Name
